# Adds carjacking data for 2022-10-30 (one more day: Oct 21 -> Oct 22 cutoff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header label to reflect the new "through" date
$ws.Name = "Through 2022-10-22"
$ws.Range("B1").Value = "October 2022 (through October 22)"

# Updated counts (existing cells incremented)
$ws.Range("B2").Value   = 3
$ws.Range("BJ2").Value  = 5
$ws.Range("L3").Value   = 6
$ws.Range("BJ4").Value  = 3
$ws.Range("L5").Value   = 10
$ws.Range("BJ6").Value  = 2
$ws.Range("B8").Value   = 3
$ws.Range("L9").Value   = 4
$ws.Range("B17").Value  = 3
$ws.Range("B43").Value  = 2
$ws.Range("L46").Value  = 4
$ws.Range("V46").Value  = 3
$ws.Range("AF47").Value = 2

# New counts (previously empty cells)
$ws.Range("AZ13").Value = 1
$ws.Range("AF16").Value = 1
$ws.Range("AP20").Value = 1
$ws.Range("BJ21").Value = 1
$ws.Range("AZ25").Value = 2
$ws.Range("V31").Value  = 1
$ws.Range("B36").Value  = 2
$ws.Range("AZ39").Value = 1
$ws.Range("AZ49").Value = 1
$ws.Range("B56").Value  = 1
$ws.Range("L77").Value  = 1
